# Insert a new data row at row 34 ("Fruta / hortaliza, semanal" update),
# shifting all existing rows from 34 downward down by one. Then populate
# the new row 34 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 34..80 down to 35..80 (inserting a blank row at 34, keeps formatting).
$ws.Rows.Item(34).Insert()

# Fill the newly inserted row 34 with the new record's data.
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C34").Value = "Ñuble"
$ws.Range("D34").Value = 44967
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100103
$ws.Range("H34").Value = "Frutos de hueso (carozo)"
$ws.Range("I34").Value = 100103002
$ws.Range("J34").Value = "Ciruela"
$ws.Range("K34").Value = "Larry Ann"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 30
$ws.Range("N34").Value = 11000
$ws.Range("O34").Value = 11000
$ws.Range("P34").Value = 11000
$ws.Range("Q34").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 611
$ws.Range("T34").Value = 18
